$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H19").Value = 1193.2222
$ws.Range("I19").Value = 1548.2
$ws.Range("J19").Value = 749.5
$ws.Range("K19").Value = 1548.2
$ws.Range("L19").Value = 749.5
$ws.Range("M19").Value = -1373.2
$ws.Range("N19").Value = -1099.5
$ws.Range("H40").Value = 8542.571
$ws.Range("I40").Value = 10159.8
$ws.Range("J40").Value = 4499.5
$ws.Range("K40").Value = 10159.8
$ws.Range("L40").Value = 4499.5
$ws.Range("M40").Value = -9984.799999999999
$ws.Range("N40").Value = -4849.5
$ws.Range("H41").Value = 729.6667
$ws.Range("I41").Value = 301
$ws.Range("K41").Value = 301
$ws.Range("M41").Value = 139
$ws.Range("H43").Value = 2270.3635
$ws.Range("I43").Value = 1486
$ws.Range("J43").Value = 2564.5
$ws.Range("K43").Value = 1486
$ws.Range("L43").Value = 2564.5
$ws.Range("M43").Value = -1417
$ws.Range("N43").Value = -2702.5
$ws.Range("H64").Value = 5359.55
$ws.Range("I64").Value = 4724.6665
$ws.Range("K64").Value = 4724.6665
$ws.Range("M64").Value = -4476.6665
$ws.Range("H67").Value = 5359.55
$ws.Range("I67").Value = 4724.6665
$ws.Range("K67").Value = 4724.6665
$ws.Range("M67").Value = -3866.6665
$ws.Range("H70").Value = 238676.92
$ws.Range("J70").Value = 100000
$ws.Range("L70").Value = 300000
$ws.Range("N70").Value = -300540
$ws.Range("H73").Value = 238676.92
$ws.Range("J73").Value = 100000
$ws.Range("L73").Value = 300000
$ws.Range("N73").Value = -301872
$ws.Range("H74").Value = 9495
$ws.Range("I74").Value = 9330
$ws.Range("K74").Value = 9330
$ws.Range("M74").Value = -8394
$ws.Range("H77").Value = 9495
$ws.Range("I77").Value = 9330
$ws.Range("K77").Value = 46650
$ws.Range("M77").Value = -41970
$ws.Range("H80").Value = 1770.875
$ws.Range("I80").Value = 1728
$ws.Range("K80").Value = 5184
$ws.Range("M80").Value = -4186
$ws.Range("H83").Value = 1770.875
$ws.Range("I83").Value = 1728
$ws.Range("K83").Value = 15552
$ws.Range("M83").Value = -10560
$ws.Range("H96").Value = 649.5
$ws.Range("I96").Value = 806.1667
$ws.Range("K96").Value = 2418.5001
$ws.Range("M96").Value = -1045.5001
$ws.Range("H98").Value = 1233.0714
$ws.Range("I98").Value = 842.2727
$ws.Range("K98").Value = 842.2727
$ws.Range("M98").Value = 655.7273
$ws.Range("H111").Value = 3428.25
$ws.Range("I111").Value = 3332.6667
$ws.Range("K111").Value = 9998.000100000001
$ws.Range("M111").Value = -6931.000100000001
$ws.Range("H113").Value = 5539.6665
$ws.Range("I113").Value = 5588.727
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 5588.727
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -2334.727
$ws.Range("N113").Value = -11508
$ws.Range("H116").Value = 3376.8333
$ws.Range("I116").Value = 3208.5
$ws.Range("K116").Value = 3208.5
$ws.Range("M116").Value = 233.5
$ws.Range("H122").Value = 1233.0714
$ws.Range("I122").Value = 842.2727
$ws.Range("K122").Value = 2526.8181
$ws.Range("M122").Value = -76.81809999999996
$ws.Range("H132").Value = 1222.6857
$ws.Range("I132").Value = 1179.5
$ws.Range("K132").Value = 3538.5
$ws.Range("M132").Value = -1008.5
$ws.Range("H135").Value = 1682.3334
$ws.Range("J135").Value = 2192.2
$ws.Range("L135").Value = 19729.8
$ws.Range("N135").Value = -24799.8

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1108.2059
$ws.Range("I2").Value = 930.4138
$ws.Range("J2").Value = 2139.4
$ws.Range("K2").Value = 930.4138
$ws.Range("L2").Value = 2139.4
$ws.Range("M2").Value = -817.4138
$ws.Range("N2").Value = -2365.4
$ws.Range("H32").Value = 9472.895
$ws.Range("I32").Value = 5589.769
$ws.Range("K32").Value = 5589.769
$ws.Range("M32").Value = -5302.769
$ws.Range("H61").Value = 4774.9287
$ws.Range("I61").Value = 4561.4
$ws.Range("J61").Value = 4893.5557
$ws.Range("K61").Value = 4561.4
$ws.Range("L61").Value = 4893.5557
$ws.Range("M61").Value = -4349.4
$ws.Range("N61").Value = -5317.5557
$ws.Range("H74").Value = 1649.6666
$ws.Range("I74").Value = 1649.6666
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1649.6666
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -775.6666
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1649.6666
$ws.Range("I77").Value = 1649.6666
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8248.333000000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3880.333000000001
$ws.Range("N77").ClearContents()
$ws.Range("H116").Value = 1108.2059
$ws.Range("I116").Value = 930.4138
$ws.Range("J116").Value = 2139.4
$ws.Range("K116").Value = 930.4138
$ws.Range("L116").Value = 2139.4
$ws.Range("M116").Value = 1363.5862
$ws.Range("N116").Value = -6727.4
$ws.Range("H122").Value = 2495.4194
$ws.Range("I122").Value = 2385.2173
$ws.Range("K122").Value = 7155.651899999999
$ws.Range("M122").Value = -4705.651899999999
$ws.Range("H130").Value = 67500
$ws.Range("J130").Value = 67500
$ws.Range("L130").Value = 67500
$ws.Range("N130").Value = -77540
$ws.Range("H132").Value = 6670.6665
$ws.Range("I132").Value = 7756
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 23268
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -20738
$ws.Range("N132").Value = -18560
$ws.Range("H136").Value = 4774.9287
$ws.Range("I136").Value = 4561.4
$ws.Range("J136").Value = 4893.5557
$ws.Range("K136").Value = 13684.2
$ws.Range("L136").Value = 14680.6671
$ws.Range("M136").Value = -11134.2
$ws.Range("N136").Value = -19780.6671

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1108.2059
$ws.Range("I3").Value = 930.4138
$ws.Range("J3").Value = 2139.4
$ws.Range("K3").Value = 930.4138
$ws.Range("L3").Value = 2139.4
$ws.Range("M3").Value = -816.4138
$ws.Range("N3").Value = -2367.4
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H86").Value = 2757.7
$ws.Range("I86").Value = 2508.5557
$ws.Range("K86").Value = 2508.5557
$ws.Range("M86").Value = -1385.5557
$ws.Range("H89").Value = 2757.7
$ws.Range("I89").Value = 2508.5557
$ws.Range("K89").Value = 12542.7785
$ws.Range("M89").Value = -6926.7785
$ws.Range("H94").Value = 2281.8
$ws.Range("J94").Value = 4386.75
$ws.Range("L94").Value = 4386.75
$ws.Range("N94").Value = -5288.75
$ws.Range("H95").Value = 11841.167
$ws.Range("J95").Value = 11841.167
$ws.Range("L95").Value = 11841.167
$ws.Range("N95").Value = -17333.167
$ws.Range("H105").Value = 4321
$ws.Range("I105").Value = 3422.9443
$ws.Range("J105").Value = 7554
$ws.Range("K105").Value = 3422.9443
$ws.Range("L105").Value = 7554
$ws.Range("M105").Value = -1675.9443
$ws.Range("N105").Value = -11048
$ws.Range("H134").Value = 5778.5
$ws.Range("I134").Value = 3963.818
$ws.Range("J134").Value = 7996.4443
$ws.Range("K134").Value = 11891.454
$ws.Range("L134").Value = 23989.3329
$ws.Range("M134").Value = -9356.454000000002
$ws.Range("N134").Value = -29059.3329

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H9").Value = 33891
$ws.Range("J9").Value = 33891
$ws.Range("L9").Value = 33891
$ws.Range("N9").Value = -34227
$ws.Range("H16").Value = 1260.8
$ws.Range("I16").Value = 1260.8
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1260.8
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -973.8
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 6012.643
$ws.Range("I31").Value = 2459.75
$ws.Range("J31").Value = 7433.8
$ws.Range("K31").Value = 2459.75
$ws.Range("L31").Value = 7433.8
$ws.Range("M31").Value = -2164.75
$ws.Range("N31").Value = -8023.8
$ws.Range("H34").Value = 6012.643
$ws.Range("I34").Value = 2459.75
$ws.Range("J34").Value = 7433.8
$ws.Range("K34").Value = 2459.75
$ws.Range("L34").Value = 7433.8
$ws.Range("M34").Value = -2257.75
$ws.Range("N34").Value = -7837.8
$ws.Range("H58").Value = 7105.1665
$ws.Range("I58").Value = 7326.3
$ws.Range("J58").Value = 5999.5
$ws.Range("K58").Value = 7326.3
$ws.Range("L58").Value = 5999.5
$ws.Range("M58").Value = -7123.3
$ws.Range("N58").Value = -6405.5
$ws.Range("H94").Value = 1288.8
$ws.Range("I94").Value = 1178.3077
$ws.Range("K94").Value = 1178.3077
$ws.Range("M94").Value = -727.3077000000001
$ws.Range("H105").Value = 4511.7144
$ws.Range("I105").Value = 4227.154
$ws.Range("K105").Value = 4227.154
$ws.Range("M105").Value = -2480.154
$ws.Range("H113").Value = 1260.8
$ws.Range("I113").Value = 1260.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1260.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 909.2
$ws.Range("N113").ClearContents()
$ws.Range("H114").Value = 37283.5
$ws.Range("J114").Value = 37283.5
$ws.Range("L114").Value = 37283.5
$ws.Range("N114").Value = -45961.5
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -52350
$ws.Range("H132").Value = 800.1667
$ws.Range("I132").Value = 800.1667
$ws.Range("K132").Value = 2400.5001
$ws.Range("M132").Value = 129.4998999999998
$ws.Range("H134").Value = 8460
$ws.Range("I134").Value = 8248.166999999999
$ws.Range("K134").Value = 24744.501
$ws.Range("M134").Value = -22209.501
$ws.Range("H136").Value = 7105.1665
$ws.Range("I136").Value = 7326.3
$ws.Range("J136").Value = 5999.5
$ws.Range("K136").Value = 21978.9
$ws.Range("L136").Value = 17998.5
$ws.Range("M136").Value = -19428.9
$ws.Range("N136").Value = -23098.5
$ws.Range("H141").Value = 411988.5
$ws.Range("J141").Value = 411988.5
$ws.Range("L141").Value = 411988.5
$ws.Range("N141").Value = -422348.5

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H14").Value = 267.58823
$ws.Range("I14").Value = 267.58823
$ws.Range("K14").Value = 802.76469
$ws.Range("M14").Value = -629.76469
$ws.Range("H128").Value = 989239.5
$ws.Range("I128").Value = 989239.5
$ws.Range("K128").Value = 2967718.5
$ws.Range("M128").Value = -2962738.5

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H59").Value = 18665.6
$ws.Range("I59").Value = 16739.334
$ws.Range("J59").Value = 21555
$ws.Range("K59").Value = 16739.334
$ws.Range("L59").Value = 21555
$ws.Range("M59").Value = -16156.334
$ws.Range("N59").Value = -22721
$ws.Range("H80").Value = 10115.5
$ws.Range("I80").Value = 6820.2
$ws.Range("J80").Value = 12650.346
$ws.Range("K80").Value = 6820.2
$ws.Range("L80").Value = 12650.346
$ws.Range("M80").Value = -5822.2
$ws.Range("N80").Value = -14646.346
$ws.Range("H83").Value = 10115.5
$ws.Range("I83").Value = 6820.2
$ws.Range("J83").Value = 12650.346
$ws.Range("K83").Value = 34101
$ws.Range("L83").Value = 63251.73
$ws.Range("M83").Value = -29109
$ws.Range("N83").Value = -73235.73
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H113").Value = 1805.1111
$ws.Range("I113").Value = 1666.3334
$ws.Range("K113").Value = 1666.3334
$ws.Range("M113").Value = 503.6666
$ws.Range("H122").Value = 113830
$ws.Range("I122").Value = 252376.5
$ws.Range("J122").Value = 2992.8
$ws.Range("K122").Value = 757129.5
$ws.Range("L122").Value = 8978.400000000001
$ws.Range("M122").Value = -754679.5
$ws.Range("N122").Value = -13878.4
$ws.Range("H132").Value = 7776.6665
$ws.Range("I132").Value = 7998.8
$ws.Range("J132").Value = 6666
$ws.Range("K132").Value = 23996.4
$ws.Range("L132").Value = 19998
$ws.Range("M132").Value = -21466.4
$ws.Range("N132").Value = -25058

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H5").Value = 17999
$ws.Range("J5").Value = 18499
$ws.Range("L5").Value = 18499
$ws.Range("N5").Value = -18725
$ws.Range("H22").Value = 999.75
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 1299.5
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 1299.5
$ws.Range("M22").Value = -405
$ws.Range("N22").Value = -1889.5
$ws.Range("H27").Value = 999.75
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 1299.5
$ws.Range("K27").Value = 700
$ws.Range("L27").Value = 1299.5
$ws.Range("M27").Value = -593
$ws.Range("N27").Value = -1513.5
$ws.Range("H31").Value = 2263.6924
$ws.Range("I31").Value = 4035.8333
$ws.Range("K31").Value = 4035.8333
$ws.Range("M31").Value = -3787.8333
$ws.Range("H40").Value = 4008.1667
$ws.Range("I40").Value = 3473.5454
$ws.Range("J40").Value = 4848.2856
$ws.Range("K40").Value = 3473.5454
$ws.Range("L40").Value = 4848.2856
$ws.Range("M40").Value = -3337.5454
$ws.Range("N40").Value = -5120.2856
$ws.Range("H46").Value = 1472.0952
$ws.Range("I46").Value = 1895.4445
$ws.Range("J46").Value = 1154.5834
$ws.Range("K46").Value = 1895.4445
$ws.Range("L46").Value = 1154.5834
$ws.Range("M46").Value = -1707.4445
$ws.Range("N46").Value = -1530.5834
$ws.Range("H122").Value = 4830.154
$ws.Range("I122").Value = 2399.5
$ws.Range("K122").Value = 7198.5
$ws.Range("M122").Value = -4748.5
$ws.Range("H132").Value = 3055.875
$ws.Range("I132").Value = 3032.1667
$ws.Range("J132").Value = 3127
$ws.Range("K132").Value = 9096.500100000001
$ws.Range("L132").Value = 9381
$ws.Range("M132").Value = -6566.500100000001
$ws.Range("N132").Value = -14441
$ws.Range("H136").Value = 6772.222
$ws.Range("I136").Value = 5850
$ws.Range("K136").Value = 17550
$ws.Range("M136").Value = -15000

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H96").Value = 3510.8333
$ws.Range("J96").Value = 3372.5
$ws.Range("L96").Value = 3372.5
$ws.Range("N96").Value = -6118.5
$ws.Range("H100").Value = 719
$ws.Range("I100").Value = 719
$ws.Range("K100").Value = 1438
$ws.Range("M100").Value = -897
$ws.Range("H113").Value = 545.625
$ws.Range("I113").Value = 606.4286
$ws.Range("J113").Value = 120
$ws.Range("K113").Value = 1819.2858
$ws.Range("L113").Value = 360
$ws.Range("M113").Value = 350.7142000000001
$ws.Range("N113").Value = -4700
$ws.Range("H122").Value = 5437.242
$ws.Range("I122").Value = 4933
$ws.Range("K122").Value = 14799
$ws.Range("M122").Value = -12349
$ws.Range("H126").Value = 1714.55
$ws.Range("I126").Value = 1716.2778
$ws.Range("J126").Value = 1699
$ws.Range("K126").Value = 5148.8334
$ws.Range("L126").Value = 5097
$ws.Range("M126").Value = -2678.8334
$ws.Range("N126").Value = -10037
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H129").Value = 40000
$ws.Range("J129").Value = 40000
$ws.Range("L129").Value = 40000
$ws.Range("N129").Value = -50000
$ws.Range("H135").Value = 59874.5
$ws.Range("J135").Value = 59874.5
$ws.Range("L135").Value = 59874.5
$ws.Range("N135").Value = -70014.5
$ws.Range("H136").Value = 5050
$ws.Range("I136").Value = 4900
$ws.Range("K136").Value = 14700
$ws.Range("M136").Value = -12150
